$d = $word.ActiveDocument

# 1. Remove the gramStart/gramEnd proofErr wrapping around "aan" by merging
#    the surrounding text into a single run of text.
$d.Content.Find.Execute(
    "worden om aan de",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "worden om aan de",
    2) | Out-Null

# 2. " in haar eerste " -> " in haar eerste, "
$d.Content.Find.Execute(
    "in haar eerste meer historische deel",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "in haar eerste, meer historische deel",
    2) | Out-Null

# 3. "gepresenteerd. Ze laat" -> "gepresenteerd; ze laat"
$d.Content.Find.Execute(
    "gepresenteerd. Ze laat",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "gepresenteerd; ze laat",
    2) | Out-Null

# 4. "behandelen en we ervoor" -> "behandelen én we ervoor"
$d.Content.Find.Execute(
    "behandelen en we ervoor",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "behandelen én we ervoor",
    2) | Out-Null
